$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so numeric-looking strings (e.g. "1.003") stay as
# literal text instead of being parsed into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.192.92'
$ws.Range("E2").Value = '  +0.05%  '
$ws.Range("D3").Value = '1.826.25'
$ws.Range("E3").Value = '  -0.04%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '234.85'
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = '0.5996'
$ws.Range("E6").Value = '  -0.75%  '
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").Value = '0.06933'
$ws.Range("E8").Value = '  -2.66%  '
$ws.Range("D9").Value = '0.2762'
$ws.Range("E9").Value = '  -1.92%  '
$ws.Range("D10").Value = '23.42'
$ws.Range("E10").Value = '  -2.45%  '
$ws.Range("D11").Value = '0.07604'
$ws.Range("E11").Value = '  -0.60%  '
$ws.Range("D12").Value = '1.847.94'
$ws.Range("E12").Value = '  +0.78%  '
$ws.Range("D13").Value = '4.719'
$ws.Range("E13").Value = '  -1.03%  '
$ws.Range("D14").Value = '0.6247'
$ws.Range("E14").Value = '  -2.47%  '
$ws.Range("D15").Value = '0.000009798'
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '77.22'
$ws.Range("E16").Value = '  -2.79%  '
$ws.Range("D17").Value = '28.963.46'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '5.541'
$ws.Range("E18").Value = '  -7.54%  '
$ws.Range("D19").Value = '215.52'
$ws.Range("E19").Value = '  -6.95%  '
$ws.Range("D20").Value = '1.004'
$ws.Range("E20").Value = '  +0.32%  '
$ws.Range("D21").Value = '11.55'
$ws.Range("E21").Value = '  -1.88%  '
$ws.Range("D22").Value = '6.835'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").Value = '1.004'
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").Value = '155.93'
$ws.Range("E24").Value = '  +0.45%  '
$ws.Range("D25").Value = '7.946'
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("D26").Value = '0.1287'
$ws.Range("E26").Value = '  +0.32%  '
$ws.Range("D27").Value = '16.47'
$ws.Range("E27").Value = '  -1.16%  '
$ws.Range("D28").Value = '0.06491'
$ws.Range("E28").Value = '  -5.14%  '
$ws.Range("D29").Value = '1.427'
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("D30").Value = '1.440'
$ws.Range("E30").Value = '  -0.69%  '
$ws.Range("D31").Value = '3.807'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").Value = '3.775'
$ws.Range("E32").Value = '  -1.72%  '
$ws.Range("D33").Value = '1.089'
$ws.Range("E33").Value = '  -4.08%  '
$ws.Range("D34").Value = '1.716'
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("D35").Value = '0.6440'
$ws.Range("E35").Value = '  -2.78%  '
$ws.Range("D36").Value = '2.542'
$ws.Range("E36").Value = '  +0.41%  '
$ws.Range("D37").Value = '2.762'
$ws.Range("E37").Value = '  +0.29%  '
$ws.Range("D38").Value = '6.596'
$ws.Range("E38").Value = '  +0.17%  '
$ws.Range("D39").Value = '0.01750'
$ws.Range("E39").Value = '  -0.64%  '
$ws.Range("D40").Value = '1.132.68'
$ws.Range("E40").Value = '  -8.27%  '
$ws.Range("D41").Value = '0.8872'
$ws.Range("E41").Value = '  -5.01%  '
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  +0.14%  '
$ws.Range("D43").Value = '1.983.32'
$ws.Range("E43").Value = '  -0.46%  '
$ws.Range("D44").Value = '100.61'
$ws.Range("E44").Value = '  +0.50%  '
$ws.Range("D45").Value = '61.88'
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").Value = '0.00000000113'
$ws.Range("E46").Value = '  -3.27%  '
$ws.Range("D47").Value = '1.609'
$ws.Range("E47").Value = '  -1.80%  '
$ws.Range("D48").Value = '8.458'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '0.05502'
$ws.Range("E49").Value = '  -1.46%  '
$ws.Range("D50").Value = '0.4528'
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").Value = '6.343'
$ws.Range("E51").Value = '  -3.54%  '

# Restore the original (default) cell style now that the text values are set,
# so no lasting formatting change is introduced.
$ws.Range("D2:E51").Style = "Normal"
